$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.9999999999999998
$ws.Range("E2").Value = 0.9999999999999998

$ws.Range("D3").Value = 0.9998659441811323
$ws.Range("E3").Value = 0.9998659441811323

$ws.Range("D4").Value = [double]"8.078271254550681E-06"
$ws.Range("E4").Value = [double]"8.078271254550681E-06"

$ws.Range("D5").Value = 0.02711021170072681
$ws.Range("E5").Value = 0.02711021170072681

$ws.Range("D6").Value = [double]"3.179700188335883E-11"
$ws.Range("E6").Value = [double]"3.179700188335883E-11"

$ws.Range("D7").Value = [double]"1.559165900070441E-05"
$ws.Range("E7").Value = 0.9999844083409993

$ws.Range("D8").Value = 0.9999999438104927
$ws.Range("E8").Value = [double]"5.618950726926641E-08"

$ws.Range("C9").Value = $false
$ws.Range("D9").Value = [double]"1.486937523797887E-06"
$ws.Range("E9").Value = 0.9999985130624762

$ws.Range("D11").Value = 0.999999999990298
$ws.Range("E11").Value = [double]"9.702016967594318E-12"
$ws.Range("F11").Value = 6.99882984161377
$ws.Range("G11").Value = 0.6
